$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.006.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.482"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.94%  "

$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.819.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.603.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.993.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "

$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.11%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  -7.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("E29").Value = "  +1.42%  "

$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("E34").Value = "  -0.72%  "

$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.83%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.27%  "

$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.493"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.05%  "

$ws.Range("E43").Value = "  +1.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.733.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.16%  "

$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.51%  "

$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0949"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.91%  "
